$d = $word.ActiveDocument

# -----------------------------------------------------------------
# Change 1: "Online Conference Room Reservation System" + "." were
# two separate runs with identical formatting; merge them into a
# single run and drop the stray _GoBack bookmark that used to sit
# at the end of that paragraph.
# -----------------------------------------------------------------
$d.Content.Find.Execute(
    "Online Conference Room Reservation System.", $true, $false, $false,
    $false, $false, $true, 1, $false,
    "Online Conference Room Reservation System.", 2) | Out-Null

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# -----------------------------------------------------------------
# Change 2: the "1.3 Definitions, Acronyms, and " heading run loses
# the comma after "Acronyms" and is split in two, with the (new)
# _GoBack bookmark landing right after "Acronyms" -- this is where
# the author's cursor ended up after editing. "Abbreviations" (the
# following run) must stay untouched/separate.
# -----------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("1.3 Definitions, Acronyms, and ", $true, $false, $false,
    $false, $false, $true, 1, $false) | Out-Null

$headingStart = $rng.Start
$acronymsEnd  = $headingStart + 25   # right after "...Acronyms"

# remove the comma that used to follow "Acronyms"
$commaRng = $d.Range($acronymsEnd, $acronymsEnd + 1)
$commaRng.Delete()

# Temporary bookmark right before "Abbreviations" so that the
# engine does not re-merge " and " with "Abbreviations" (both runs
# share identical formatting and would otherwise collapse back into
# one run).
$beforeAbbrevPos = $acronymsEnd + 5   # length of " and "
$tempRng = $d.Range($beforeAbbrevPos, $beforeAbbrevPos)
$d.Bookmarks.Add("ZZTempSplit", $tempRng) | Out-Null

# Place the real _GoBack bookmark right after "Acronyms" (collapsed,
# splitting the run there).
$gobackRng = $d.Range($acronymsEnd, $acronymsEnd)
$d.Bookmarks.Add("_GoBack", $gobackRng) | Out-Null

# Drop the scaffolding bookmark; the run split it enforced remains.
if ($d.Bookmarks.Exists("ZZTempSplit")) {
    $d.Bookmarks("ZZTempSplit").Delete()
}
